$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in homework scores for rows 21 and 22
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 5

$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 5
$ws.Range("F22").Value = 5

# Update the active selection to E21 to match the saved view state
$ws.Range("E21").Select()
